$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Insert two new rows at row 26 (pushes former rows 26-29 down to 28-31)
#    to make room for the new "ifExpr" test rows.
# ---------------------------------------------------------------------------
$ws.Rows(26).Resize(2).Insert()

# Apply the same formatting used by the other "ifTest"-style rows (row 3) to
# the two freshly inserted rows so every column picks up the bordered style.
$ws.Range("A3:H3").Copy()
$ws.Range("A26:H27").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Fill in the new ifExpr rows.
#    Values are written in a specific order so that new shared strings are
#    created in the same order as the target workbook:
#      ifExpr, apple, vars["apple"] == 3, runLoopEnd,
#      Test RunLoop start end step, Test If Expr positive, Test If Expr negative
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = "ifExpr"
$ws.Range("F26").Value = "apple"
$ws.Range("H26").Value = 'vars["apple"] == 3'
$ws.Range("A26").Value = 25
$ws.Range("D26").Value = "yes"
$ws.Range("G26").Value = 3

$ws.Range("A27").Value = 26
$ws.Range("C27").Value = "ifExpr"
$ws.Range("D27").Value = "no"
$ws.Range("F27").Value = "apple"
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 'vars["apple"] == 3'

# ---------------------------------------------------------------------------
# 3. Append six new "runLoopEnd" rows after the (now shifted) row 31.
# ---------------------------------------------------------------------------
$ws.Range("A28:H28").Copy()
$ws.Range("A32:H37").PasteSpecial(-4122)

$ws.Range("C32").Value = "runLoopEnd"
$ws.Range("B32").Value = "Test RunLoop start end step"

$ws.Range("A32").Value = 31
$ws.Range("D32").Value = 5
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 5
$ws.Range("H32").Value = 1

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "Test RunLoop start end step"
$ws.Range("C33").Value = "runLoopEnd"
$ws.Range("D33").Value = 3
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 5
$ws.Range("H33").Value = 2

$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "Test RunLoop start end step"
$ws.Range("C34").Value = "runLoopEnd"
$ws.Range("D34").Value = 2
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 5
$ws.Range("H34").Value = 3

$ws.Range("A35").Value = 34
$ws.Range("B35").Value = "Test RunLoop start end step"
$ws.Range("C35").Value = "runLoopEnd"
$ws.Range("D35").Value = 2
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 5
$ws.Range("H35").Value = 4

$ws.Range("A36").Value = 35
$ws.Range("B36").Value = "Test RunLoop start end step"
$ws.Range("C36").Value = "runLoopEnd"
$ws.Range("D36").Value = 1
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 5

$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "Test RunLoop start end step"
$ws.Range("C37").Value = "runLoopEnd"
$ws.Range("D37").Value = 1
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 5
$ws.Range("H37").Value = 6

# ---------------------------------------------------------------------------
# 4. Finally add the two "Test Name" values for the ifExpr rows, after the
#    runLoopEnd block, to match the shared-string ordering of the target file.
# ---------------------------------------------------------------------------
$ws.Range("B26").Value = "Test If Expr positive"
$ws.Range("B27").Value = "Test If Expr negative"

# ---------------------------------------------------------------------------
# 5. Update the sheet view: scroll position and active selection.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("B18").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
